$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5000
$ws.Range("I32").Value = 5000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4674
$ws.Range("N32").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 26571.428
$ws.Range("J44").Value = 26571.428
$ws.Range("L44").Value = 26571.428
$ws.Range("N44").Value = -27495.428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 15386250
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = 0

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 5537.923
$ws.Range("J113").Value = 5863.364
$ws.Range("L113").Value = 5863.364
$ws.Range("N113").Value = -12371.364

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 299982.06
$ws.Range("I116").Value = 669436.25
$ws.Range("J116").Value = 8307.684999999999
$ws.Range("K116").Value = 669436.25
$ws.Range("L116").Value = 8307.684999999999
$ws.Range("M116").Value = -665994.25
$ws.Range("N116").Value = -15191.685

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 1430.75
$ws.Range("I127").Value = 698.375
$ws.Range("J127").Value = 1796.9375
$ws.Range("K127").Value = 2095.125
$ws.Range("L127").Value = 5390.8125
$ws.Range("M127").Value = 2864.875
$ws.Range("N127").Value = -15310.8125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2466.926
$ws.Range("I137").Value = 1713.2222
$ws.Range("K137").Value = 5139.6666
$ws.Range("M137").Value = -2589.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3415.6538
$ws.Range("I32").Value = 3240.41
$ws.Range("J32").Value = 4044.4707
$ws.Range("K32").Value = 3240.41
$ws.Range("L32").Value = 4044.4707
$ws.Range("M32").Value = -2953.41
$ws.Range("N32").Value = -4618.4707

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3810.2646
$ws.Range("I74").Value = 4576.95
$ws.Range("J74").Value = 2715
$ws.Range("K74").Value = 4576.95
$ws.Range("L74").Value = 2715
$ws.Range("M74").Value = -3702.95
$ws.Range("N74").Value = -4463

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3810.2646
$ws.Range("I77").Value = 4576.95
$ws.Range("J77").Value = 2715
$ws.Range("K77").Value = 22884.75
$ws.Range("L77").Value = 13575
$ws.Range("M77").Value = -18516.75
$ws.Range("N77").Value = -22311

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1815.7142
$ws.Range("I132").Value = 1023.4722
$ws.Range("K132").Value = 3070.4166
$ws.Range("M132").Value = -540.4166

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 19418.125
$ws.Range("J81").Value = 19418.125
$ws.Range("L81").Value = 19418.125
$ws.Range("N81").Value = -21540.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H84").Value = 19418.125
$ws.Range("J84").Value = 19418.125
$ws.Range("L84").Value = 58254.375
$ws.Range("N84").Value = -68862.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3031.75
$ws.Range("I31").Value = 1123.8182
$ws.Range("J31").Value = 4646.154
$ws.Range("K31").Value = 1123.8182
$ws.Range("L31").Value = 4646.154
$ws.Range("M31").Value = -828.8181999999999
$ws.Range("N31").Value = -5236.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3031.75
$ws.Range("I34").Value = 1123.8182
$ws.Range("J34").Value = 4646.154
$ws.Range("K34").Value = 1123.8182
$ws.Range("L34").Value = 4646.154
$ws.Range("M34").Value = -921.8181999999999
$ws.Range("N34").Value = -5050.154

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1419.7234
$ws.Range("I58").Value = 1357.4304
$ws.Range("J58").Value = 1747.8
$ws.Range("K58").Value = 1357.4304
$ws.Range("L58").Value = 1747.8
$ws.Range("M58").Value = -1154.4304
$ws.Range("N58").Value = -2153.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 670.3182
$ws.Range("I107").Value = 565.6875
$ws.Range("J107").Value = 949.3333
$ws.Range("K107").Value = 565.6875
$ws.Range("L107").Value = 949.3333
$ws.Range("M107").Value = 1354.3125
$ws.Range("N107").Value = -4789.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2770.2727
$ws.Range("I132").Value = 1780.5
$ws.Range("J132").Value = 7224.25
$ws.Range("K132").Value = 5341.5
$ws.Range("L132").Value = 21672.75
$ws.Range("M132").Value = -2811.5
$ws.Range("N132").Value = -26732.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1419.7234
$ws.Range("I136").Value = 1357.4304
$ws.Range("J136").Value = 1747.8
$ws.Range("K136").Value = 4072.2912
$ws.Range("L136").Value = 5243.4
$ws.Range("M136").Value = -1522.2912
$ws.Range("N136").Value = -10343.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 73782.5
$ws.Range("J140").Value = 73782.5
$ws.Range("L140").Value = 73782.5
$ws.Range("N140").Value = -84142.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 769.1613
$ws.Range("I113").Value = 763.76
$ws.Range("J113").Value = 791.6667
$ws.Range("K113").Value = 2291.28
$ws.Range("L113").Value = 2375.0001
$ws.Range("M113").Value = -121.2799999999997
$ws.Range("N113").Value = -6715.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2248.803
$ws.Range("J121").Value = 2312.8594
$ws.Range("L121").Value = 6938.5782
$ws.Range("N121").Value = -9558.5782

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5102848.5
$ws.Range("I131").Value = 125000250
$ws.Range("J131").Value = 831.34045
$ws.Range("K131").Value = 375000750
$ws.Range("L131").Value = 2494.02135
$ws.Range("M131").Value = -374995710
$ws.Range("N131").Value = -12574.02135

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1909.3125
$ws.Range("I139").Value = 1186.2727
$ws.Range("K139").Value = 3558.8181
$ws.Range("M139").Value = 1581.1819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 29000
$ws.Range("J4").Value = 29000
$ws.Range("L4").Value = 29000
$ws.Range("N4").Value = -29224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 12859829
$ws.Range("I10").Value = 15000000
$ws.Range("K10").Value = 15000000
$ws.Range("M10").Value = -14999831

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 758
$ws.Range("I97").Value = 757.0909
$ws.Range("K97").Value = 757.0909
$ws.Range("M97").Value = -261.0909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 60865
$ws.Range("I22").Value = 112466.89
$ws.Range("J22").Value = 2812.875
$ws.Range("K22").Value = 112466.89
$ws.Range("L22").Value = 2812.875
$ws.Range("M22").Value = -112171.89
$ws.Range("N22").Value = -3402.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 60865
$ws.Range("I27").Value = 112466.89
$ws.Range("J27").Value = 2812.875
$ws.Range("K27").Value = 112466.89
$ws.Range("L27").Value = 2812.875
$ws.Range("M27").Value = -112359.89
$ws.Range("N27").Value = -3026.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1255.875
$ws.Range("I61").Value = 1357
$ws.Range("J61").Value = 952.5
$ws.Range("K61").Value = 1357
$ws.Range("L61").Value = 952.5
$ws.Range("M61").Value = -1155
$ws.Range("N61").Value = -1356.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1255.875
$ws.Range("I113").Value = 1357
$ws.Range("J113").Value = 952.5
$ws.Range("K113").Value = 1357
$ws.Range("L113").Value = 952.5
$ws.Range("M113").Value = 813
$ws.Range("N113").Value = -5292.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3188.1667
$ws.Range("I136").Value = 1274.9412
$ws.Range("K136").Value = 3824.8236
$ws.Range("M136").Value = -1274.8236

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 56674
$ws.Range("J29").Value = 70011
$ws.Range("L29").Value = 70011
$ws.Range("N29").Value = -70591

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11906545
$ws.Range("I132").Value = 1039.5217
$ws.Range("J132").Value = 66671868
$ws.Range("K132").Value = 3118.5651
$ws.Range("L132").Value = 200015604
$ws.Range("M132").Value = -588.5650999999998
$ws.Range("N132").Value = -200020664

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 35915.832
$ws.Range("J133").Value = 35915.832
$ws.Range("L133").Value = 35915.832
$ws.Range("N133").Value = -46035.832
